$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8 for "un_franzosa_ControlvsCD_ConvCD"
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "un_franzosa_ControlvsCD_ConvCD"
$ws.Range("B8").Value = 0.06
$ws.Range("C8").Value = 0.04
$ws.Range("D8").Value = 0.02
$ws.Range("E8").Value = 0.38
$ws.Range("F8").Value = 0.88
$ws.Range("G8").Value = 0.52
$ws.Range("H8").Value = 0.54

# Insert a new row at row 13 for "un_franzosa_ControlvsUC_ConvUC"
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = "un_franzosa_ControlvsUC_ConvUC"
$ws.Range("B13").Value = 0.02
$ws.Range("C13").Value = 0.06
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0.58
$ws.Range("F13").Value = 0.92
$ws.Range("G13").Value = 0.34
$ws.Range("H13").Value = 0.4
